$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 40; existing rows 40-169 shift down to 41-170.
$ws.Rows("40:40").Insert()

# Populate the new row 40 with the new weekly record.
$ws.Range("A40").Value = 7
$ws.Range("B40").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C40").Value = "Ñuble"
$ws.Range("D40").Value = 44487
$ws.Range("E40").Value = 16
$ws.Range("F40").Value = 100112008
$ws.Range("G40").Value = "Coliflor"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 700
$ws.Range("L40").Value = 800
$ws.Range("M40").Value = 750
$ws.Range("N40").Value = "$/unidad"
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 750
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"
